$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for every data row
# (rows 2-176). The source data was refreshed, bumping this date from
# 2023-09-01 (serial 45170) to 2023-09-05 (serial 45174) for every row.
$newDate = Get-Date -Year 2023 -Month 9 -Day 5 -Hour 0 -Minute 0 -Second 0

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 176 }

$ws.Range("C2:C$lastRow").Value = $newDate
